$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 77; $r++) {
    $ws.Cells.Item($r, 3).Value2 = 45181
}
